$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For each updated cell, force a text number format before assigning the new value so that
# numeric-looking text (e.g. "309.57", "-2.01%") is stored as text, matching the original
# inline-string cells, instead of being auto-converted to a number/percentage by Excel.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '309.57'
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '-2.01%'
$c = $ws.Range("G2")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '37.70'
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '-4.33%'
$c = $ws.Range("G3")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '5.063'
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '-1.42%'
$c = $ws.Range("G4")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '0.07771'
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '-4.91%'
$c = $ws.Range("G5")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '4.349'
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '-0.64%'
$c = $ws.Range("G6")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '1.897'
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '-3.99%'
$c = $ws.Range("G7")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '8.184'
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '-1.85%'
$c = $ws.Range("G8")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '3.110'
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '-6.31%'
$c = $ws.Range("G9")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.9206'
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '-1.69%'
$c = $ws.Range("G10")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.1254'
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '-3.93%'
$c = $ws.Range("G11")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.1882'
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '-4.59%'
$c = $ws.Range("G12")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.08815'
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '-2.59%'
$c = $ws.Range("G13")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.03412'
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '-2.45%'
$c = $ws.Range("G14")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.09717'
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '-0.49%'
$c = $ws.Range("G15")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.001366'
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '-3.33%'
$c = $ws.Range("G16")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.006072'
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '0.45%'
$c = $ws.Range("G17")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '3.596'
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '-1.40%'
$c = $ws.Range("G18")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.3413'
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '-2.22%'
$c = $ws.Range("G19")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '5.023'
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '0.36%'
$c = $ws.Range("G20")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("G21")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.2587'
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '4.26%'
$c = $ws.Range("G22")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.02102'
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '5,583.31%'
$c = $ws.Range("G23")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '0.04383'
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '0.57%'
$c = $ws.Range("G24")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '0.001210'
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '-2.81%'
$c = $ws.Range("G25")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '0.004266'
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '-10.49%'
$c = $ws.Range("G26")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '0.0001349'
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '3.62%'
$c = $ws.Range("G27")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("G28")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("G29")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("G30")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("G31")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("G32")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("G33")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("G34")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("G35")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("G36")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("G37")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("G38")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.02160'
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '-3.50%'
$c = $ws.Range("G39")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.04971'
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '-3.76%'
$c = $ws.Range("G40")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.008050'
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '3.67%'
$c = $ws.Range("G41")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.009927'
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '-4.95%'
$c = $ws.Range("G42")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.1345'
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '-3.97%'
$c = $ws.Range("G43")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.002058'
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '0.67%'
$c = $ws.Range("G44")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.008689'
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '-6.35%'
$c = $ws.Range("G45")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.00006431'
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '-7.60%'
$c = $ws.Range("G46")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.00000000749'
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '-0.25%'
$c = $ws.Range("G47")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.003368'
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '16.71%'
$c = $ws.Range("G48")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.001689'
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '-0.26%'
$c = $ws.Range("G49")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.00002098'
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '-0.25%'
$c = $ws.Range("G50")
$c.NumberFormat = "@"
$c.Value = '2'

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.0001999'
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '-0.25%'
$c = $ws.Range("G51")
$c.NumberFormat = "@"
$c.Value = '2'
